$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the discount value (row 18 "Client discount", column E = UNIT PRICE/AMOUNT)
$ws.Range("E18").Value = 125

# Write the signature (name, email address) into the bottom signature cell
$ws.Range("A31").Value = "RPA Dev, developer.rpa@mail.com"
